# Add a new work-log entry (row 49) to Sheet1 and move the "latest entry"
# formatting down from the old last row (48) to the new one (49), mirroring
# what Excel does when a user types a new row right after the previous
# last row of a manually-maintained table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1) Carry the existing (old) last-row formatting for columns A:C down
#        onto the brand-new row 49, before we touch row 48's own formatting.
$ws.Range("A48:C48").Copy()
$ws.Range("A49:C49").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- 2) Revert the old last row (48, columns B:C) back to the plain style
#        used by every earlier row (it is no longer the "latest" entry).
$ws.Range("B47:C47").Copy()
$ws.Range("B48:C48").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- 3) Enter the new log entry's data.
$ws.Cells.Item(49, 1).Value = 45755
$ws.Cells.Item(49, 2).Value = 4
$ws.Cells.Item(49, 3).Value = "Debugging minor errors, generated evaluation form"

# --- 4) Match the saved cursor position left behind in the sheet.
$ws.Range("C53").Select() | Out-Null
